$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update confidentiality footer date (2021-04-05 -> 2021-04-06)
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-06 for illustrative purposes only and are subject to change."

# Update Weight (D) / Percent Change (E) columns for holdings rows 2-34
$ws.Range("D2").Value = 0.03838588865667181
$ws.Range("E2").Value = 0.001542614731970593
$ws.Range("D3").Value = 0.02173327918603344
$ws.Range("E3").Value = 0.002318392581143858
$ws.Range("D4").Value = 0.01999775199704056
$ws.Range("E4").Value = 0.009410801963993354
$ws.Range("D5").Value = 0.03967919016770016
$ws.Range("E5").Value = 0
$ws.Range("D6").Value = 0.03712843891697935
$ws.Range("E6").Value = 0.001950078003120304
$ws.Range("D7").Value = 0.02082874868489674
$ws.Range("E7").Value = 0.004264392324093702
$ws.Range("D8").Value = 0.03826503696500651
$ws.Range("E8").Value = -0.004361098996947366
$ws.Range("D9").Value = 0.0211738623576211
$ws.Range("E9").Value = 0.003386417719201784
$ws.Range("D10").Value = 0.02608943089781018
$ws.Range("E10").Value = 0.006602702981532982
$ws.Range("D11").Value = 0.0243251038624739
$ws.Range("E11").Value = -0.001106500691562928
$ws.Range("D12").Value = 0.05888621143027681
$ws.Range("E12").Value = -0.002196729314132173
$ws.Range("D13").Value = 0.02654969012441964
$ws.Range("E13").Value = 0.001094890510948909
$ws.Range("D14").Value = 0.02765926476522365
$ws.Range("E14").Value = 0.007181615065432467
$ws.Range("D15").Value = 0.03549346894849015
$ws.Range("E15").Value = 0.004076779344317938
$ws.Range("D16").Value = 0.01894416210426858
$ws.Range("E16").Value = -0.006339581036383768
$ws.Range("D17").Value = 0.02981995305032899
$ws.Range("E17").Value = 0.004603303547251469
$ws.Range("D18").Value = 0.02388503144404906
$ws.Range("E18").Value = 0.004597701149425371
$ws.Range("D19").Value = 0.1332940865894565
$ws.Range("E19").Value = 0.001326259946949682
$ws.Range("D20").Value = 0.009533556611992535
$ws.Range("E20").Value = -0.01571428571428568
$ws.Range("D21").Value = 0.01603446248828088
$ws.Range("E21").Value = 0.004230118443316444
$ws.Range("D22").Value = 0.01702060152597197
$ws.Range("E22").Value = 0.002742082907683185
$ws.Range("D23").Value = 0.0164880466103173
$ws.Range("E23").Value = -0.001067615658362997
$ws.Range("D24").Value = 0.02174711387857597
$ws.Range("E24").Value = 0.002660989888238374
$ws.Range("D25").Value = 0.01209690443102054
$ws.Range("E25").Value = 0.005963029218843285
$ws.Range("D26").Value = 0.04356959184104121
$ws.Range("E26").Value = 0.003391526743022233
$ws.Range("D27").Value = 0.02522554317048519
$ws.Range("E27").Value = 0.0001963286541668552
$ws.Range("D28").Value = 0.04797419189246105
$ws.Range("E28").Value = 0.002868068833651982
$ws.Range("D29").Value = 0.05889552427778598
$ws.Range("E29").Value = -0.006986743102830428
$ws.Range("D30").Value = 0.01331468036369415
$ws.Range("E30").Value = 0.002013422818792021
$ws.Range("D31").Value = 0.01470272529446901
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0.04408411320805095
$ws.Range("E32").Value = 0.002582644628099207
$ws.Range("D33").Value = 0.01717434425710608
$ws.Range("E33").Value = -0.009061588081707872
$ws.Range("E34").Value = 0.0009969591668781419

$ws.Protect()
